$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.580986499786377
$ws.Range("B1").Value = 2.763327836990356
$ws.Range("C1").Value = 2.371453046798706
$ws.Range("D1").Value = 2.514743328094482
$ws.Range("E1").Value = 2.786831140518188
